$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 944.9
$ws.Range("I6").Value = 109.583336
$ws.Range("J6").Value = 3092.8572
$ws.Range("K6").Value = 328.750008
$ws.Range("L6").Value = 9278.571599999999
$ws.Range("M6").Value = -216.750008
$ws.Range("N6").Value = -9502.571599999999
$ws.Range("H70").Value = 1558.2106
$ws.Range("I70").Value = 1386.4615
$ws.Range("J70").Value = 1930.3334
$ws.Range("K70").Value = 4159.3845
$ws.Range("L70").Value = 5791.0002
$ws.Range("M70").Value = -3889.3845
$ws.Range("N70").Value = -6331.0002
$ws.Range("H73").Value = 1558.2106
$ws.Range("I73").Value = 1386.4615
$ws.Range("J73").Value = 1930.3334
$ws.Range("K73").Value = 4159.3845
$ws.Range("L73").Value = 5791.0002
$ws.Range("M73").Value = -3223.3845
$ws.Range("N73").Value = -7663.0002
$ws.Range("H82").Value = 952.36365
$ws.Range("I82").Value = 952.36365
$ws.Range("K82").Value = 2857.09095
$ws.Range("M82").Value = -2451.09095
$ws.Range("H85").Value = 952.36365
$ws.Range("I85").Value = 952.36365
$ws.Range("K85").Value = 2857.09095
$ws.Range("M85").Value = -1453.09095
$ws.Range("H115").Value = 3554.8333
$ws.Range("I115").Value = 2305.8
$ws.Range("K115").Value = 6917.400000000001
$ws.Range("M115").Value = -5350.400000000001
$ws.Range("H132").Value = 5793.8696
$ws.Range("I132").Value = 5393.381
$ws.Range("K132").Value = 16180.143
$ws.Range("M132").Value = -13650.143
$ws.Range("H137").Value = 4337.1304
$ws.Range("I137").Value = 3967.9412
$ws.Range("J137").Value = 5383.1665
$ws.Range("K137").Value = 11903.8236
$ws.Range("L137").Value = 16149.4995
$ws.Range("M137").Value = -9353.8236
$ws.Range("N137").Value = -21249.4995
$ws.Range("H138").Value = 1915.07
$ws.Range("I138").Value = 696.1429000000001
$ws.Range("J138").Value = 2113.5
$ws.Range("K138").Value = 2088.4287
$ws.Range("L138").Value = 6340.5
$ws.Range("M138").Value = 3051.5713
$ws.Range("N138").Value = -16620.5
$ws.Range("H139").Value = 45926.668
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 45926.668
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = ""
$ws.Range("M139").Value = ""
$ws.Range("N139").Value = -56206.668
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 681755.0600000001
$ws.Range("I32").Value = 813303.25
$ws.Range("J32").Value = 24014.25
$ws.Range("K32").Value = 813303.25
$ws.Range("L32").Value = 24014.25
$ws.Range("M32").Value = -813016.25
$ws.Range("N32").Value = -24588.25
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = ""
$ws.Range("N64").Value = ""
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = ""
$ws.Range("N67").Value = ""
$ws.Range("H74").Value = 1791.0416
$ws.Range("J74").Value = 2666.3333
$ws.Range("L74").Value = 2666.3333
$ws.Range("N74").Value = -4414.3333
$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100676
$ws.Range("H77").Value = 1791.0416
$ws.Range("J77").Value = 2666.3333
$ws.Range("L77").Value = 13331.6665
$ws.Range("N77").Value = -22067.6665
$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102340
$ws.Range("H132").Value = 3238.487
$ws.Range("I132").Value = 2917.4644
$ws.Range("J132").Value = 4055.6365
$ws.Range("K132").Value = 8752.393199999999
$ws.Range("L132").Value = 12166.9095
$ws.Range("M132").Value = -6222.393199999999
$ws.Range("N132").Value = -17226.9095
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2990.8
$ws.Range("I20").Value = 2756
$ws.Range("J20").Value = 3277.7778
$ws.Range("K20").Value = 2756
$ws.Range("L20").Value = 3277.7778
$ws.Range("M20").Value = -2509
$ws.Range("N20").Value = -3771.7778
$ws.Range("H134").Value = 2117.6052
$ws.Range("I134").Value = 1956.2572
$ws.Range("K134").Value = 5868.7716
$ws.Range("M134").Value = -3333.7716
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5693.9434
$ws.Range("I31").Value = 1502.7241
$ws.Range("J31").Value = 10758.333
$ws.Range("K31").Value = 1502.7241
$ws.Range("L31").Value = 10758.333
$ws.Range("M31").Value = -1207.7241
$ws.Range("N31").Value = -11348.333
$ws.Range("H34").Value = 5693.9434
$ws.Range("I34").Value = 1502.7241
$ws.Range("J34").Value = 10758.333
$ws.Range("K34").Value = 1502.7241
$ws.Range("L34").Value = 10758.333
$ws.Range("M34").Value = -1300.7241
$ws.Range("N34").Value = -11162.333
$ws.Range("H107").Value = 1603187.6
$ws.Range("I107").Value = 4167312.2
$ws.Range("J107").Value = 609.75
$ws.Range("K107").Value = 4167312.2
$ws.Range("L107").Value = 609.75
$ws.Range("M107").Value = -4165392.2
$ws.Range("N107").Value = -4449.75
$ws.Range("H134").Value = 1912.4762
$ws.Range("I134").Value = 1703.4286
$ws.Range("J134").Value = 2330.5715
$ws.Range("K134").Value = 5110.2858
$ws.Range("L134").Value = 6991.7145
$ws.Range("M134").Value = -2575.2858
$ws.Range("N134").Value = -12061.7145
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 668.5
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 3000
$ws.Range("N5").Value = -3224
$ws.Range("H39").Value = 1490.279
$ws.Range("J39").Value = 1490.279
$ws.Range("L39").Value = 4470.837
$ws.Range("N39").Value = -5058.837
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = ""
$ws.Range("N98").Value = ""
$ws.Range("H122").Value = 14581.143
$ws.Range("I122").Value = 344.83334
$ws.Range("J122").Value = 99999
$ws.Range("K122").Value = 3103.50006
$ws.Range("L122").Value = 899991
$ws.Range("M122").Value = -653.5000600000003
$ws.Range("N122").Value = -904891
$ws.Range("H135").Value = 668.5
$ws.Range("J135").Value = 1000
$ws.Range("L135").Value = 9000
$ws.Range("N135").Value = -14070
$ws.Range("H137").Value = 9486.733
$ws.Range("J137").Value = 3983.25
$ws.Range("L137").Value = 11949.75
$ws.Range("N137").Value = -22149.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4000
$ws.Range("I5").Value = 2000
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -1888
$ws.Range("N5").Value = -5224
$ws.Range("H46").Value = 4192.381
$ws.Range("J46").Value = 4349.95
$ws.Range("L46").Value = 4349.95
$ws.Range("N46").Value = -4661.95
$ws.Range("H132").Value = 3572.75
$ws.Range("I132").Value = 3286.5454
$ws.Range("K132").Value = 9859.636200000001
$ws.Range("M132").Value = -7329.636200000001
$ws.Range("H138").Value = 69321
$ws.Range("J138").Value = 69321
$ws.Range("L138").Value = 69321
$ws.Range("N138").Value = -79601
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2781
$ws.Range("I7").Value = 1966.6666
$ws.Range("J7").Value = 4002.5
$ws.Range("K7").Value = 1966.6666
$ws.Range("L7").Value = 4002.5
$ws.Range("M7").Value = -1854.6666
$ws.Range("N7").Value = -4226.5
$ws.Range("H126").Value = 2781
$ws.Range("I126").Value = 1966.6666
$ws.Range("J126").Value = 4002.5
$ws.Range("K126").Value = 5899.9998
$ws.Range("L126").Value = 12007.5
$ws.Range("M126").Value = -3429.9998
$ws.Range("N126").Value = -16947.5
$ws.Range("H132").Value = 3167.2917
$ws.Range("I132").Value = 2026.375
$ws.Range("J132").Value = 5449.125
$ws.Range("K132").Value = 6079.125
$ws.Range("L132").Value = 16347.375
$ws.Range("M132").Value = -3549.125
$ws.Range("N132").Value = -21407.375
$ws.Range("H138").Value = 49800
$ws.Range("J138").Value = 49800
$ws.Range("L138").Value = 49800
$ws.Range("N138").Value = -60080
$ws.Range("H139").Value = 4870479
$ws.Range("J139").Value = 60638.332
$ws.Range("L139").Value = 60638.332
$ws.Range("N139").Value = -70918.33199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 59356
$ws.Range("J16").Value = 59356
$ws.Range("L16").Value = 59356
$ws.Range("N16").Value = -59940
$ws.Range("H122").Value = 3025.1365
$ws.Range("I122").Value = 1776.04
$ws.Range("J122").Value = 4668.684
$ws.Range("K122").Value = 5328.12
$ws.Range("L122").Value = 14006.052
$ws.Range("M122").Value = -2878.12
$ws.Range("N122").Value = -18906.052
$ws.Range("H136").Value = 2196.258
$ws.Range("I136").Value = 1964.7693
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 5894.3079
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -3344.3079
$ws.Range("N136").Value = -15300
